$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.517.23"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.624.39"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.59"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0610"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.22"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "1.852.48"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "1.617.21"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.512"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("D17").Value = "26.510.08"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.58"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.76"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.32"
$ws.Range("E22").Value = "  -1.60%  "
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.16"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.99"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.07"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.65"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0494"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").Value = "1.524.42"
$ws.Range("E32").Value = "  +4.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.26"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.01"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.568"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").Value = "1.764.45"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.79"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.761"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("E46").Value = "  -5.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.89"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0502"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0964"
$ws.Range("E51").Value = "  -0.03%  "
